{"js": "// Underline the primary-key field name(s) called out (via highlight colour)\n// in the short \"TABLE(Field1, Field2, ...)\" schema-summary lines at the top\n// of the document. Each of the first seven paragraphs introduces one table;\n// the run(s) that spell out its primary key (already highlighted in a\n// distinct colour) should additionally become underlined. For the two\n// paragraphs whose primary key is a composite of two highlighted fields\n// (CHITIETKHO / CHITIETHOADON), the \", \" separator that sits between the two\n// highlighted fields is underlined as well.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// paragraph index (0-based) -> exact text span to underline within that\n// paragraph. The span is searched for within the paragraph's own range, so\n// it is safe even though some of these substrings recur elsewhere in the\n// document (e.g. later in the descriptive paragraphs).\nconst targets = [\n  { paragraph: 0, text: \"MaSP\" }, // SANPHAM(MaSP, ...)\n  { paragraph: 1, text: \"MaNV\" }, // NHANVIEN(MaNV, ...)\n  { paragraph: 2, text: \"SDTKH\" }, // KHACHHANG(SDTKH, ...)\n  { paragraph: 3, text: \"MaHD\" }, // HOADON(MaHD, ...)\n  { paragraph: 4, text: \"MaLoHang\" }, // KHOVAN(MaLoHang, ...)\n  { paragraph: 5, text: \"MaLoHang, MaSP\" }, // CHITIETKHO(MaLoHang, MaSP, ...)\n  { paragraph: 6, text: \"MaHD, MaSP\" }, // CHITIETHOADON(MaHD, MaSP, ...)\n];\n\nconst foundRanges = [];\nfor (const target of targets) {\n  const paragraph = paragraphs.items[target.paragraph];\n  const paragraphRange = paragraph.getRange();\n  const searchResults = paragraphRange.search(target.text, { matchCase: true });\n  searchResults.load(\"items\");\n  foundRanges.push(searchResults);\n}\nawait context.sync();\n\nfor (const searchResults of foundRanges) {\n  if (searchResults.items.length > 0) {\n    searchResults.items[0].font.underline = Word.UnderlineType.single;\n  }\n}\nawait context.sync();\n", "ps1": "# Underline the primary-key field name(s) called out (via highlight colour)\n# in the short \"TABLE(Field1, Field2, ...)\" schema-summary lines at the top\n# of the document. Each of the first seven paragraphs introduces one table;\n# the run(s) that spell out its primary key (already highlighted in a\n# distinct colour) should additionally become underlined. For the two\n# paragraphs whose primary key is a composite of two highlighted fields\n# (CHITIETKHO / CHITIETHOADON), the \", \" separator that sits between the two\n# highlighted fields is underlined as well.\n\n$d = $word.ActiveDocument\n\n# Paragraph number (1-based, matching $d.Paragraphs.Item(N)) -> exact text\n# span to underline within that paragraph. The span is located with a Find\n# scoped to the paragraph's own Range, so it is safe even though some of\n# these substrings recur elsewhere in the document (e.g. later in the\n# descriptive paragraphs).\n$targets = @(\n    @{ Paragraph = 1; Text = \"MaSP\" },            # SANPHAM(MaSP, ...)\n    @{ Paragraph = 2; Text = \"MaNV\" },             # NHANVIEN(MaNV, ...)\n    @{ Paragraph = 3; Text = \"SDTKH\" },            # KHACHHANG(SDTKH, ...)\n    @{ Paragraph = 4; Text = \"MaHD\" },             # HOADON(MaHD, ...)\n    @{ Paragraph = 5; Text = \"MaLoHang\" },         # KHOVAN(MaLoHang, ...)\n    @{ Paragraph = 6; Text = \"MaLoHang, MaSP\" },   # CHITIETKHO(MaLoHang, MaSP, ...)\n    @{ Paragraph = 7; Text = \"MaHD, MaSP\" }        # CHITIETHOADON(MaHD, MaSP, ...)\n)\n\nforeach ($target in $targets) {\n    $p = $d.Paragraphs.Item($target.Paragraph)\n    $r = $p.Range.Duplicate()\n    $r.Find.ClearFormatting()\n    $r.Find.Text = $target.Text\n    $r.Find.Forward = $true\n    $r.Find.Wrap = 0\n    $found = $r.Find.Execute()\n    if ($found) {\n        $r.Font.Underline = 1\n    }\n}\n"}
